$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.719134520185651
$ws.Range("C2").Value = 4.341831986390855
$ws.Range("E2").Value = 19.15624232107774
$ws.Range("F2").Value = 42.43582529800561
$ws.Range("G2").Value = 37.37058304504887
$ws.Range("H2").Value = 16.46480571249877
$ws.Range("I2").Value = 24.12194164514337
$ws.Range("J2").Value = 8.690736112842444
$ws.Range("K2").Value = 9.055337863823665
$ws.Range("M2").Value = 17.46026666267149

$ws.Range("B3").Value = 8.465870614906843
$ws.Range("C3").Value = 4.15749407233841
$ws.Range("E3").Value = 19.06452695066666
$ws.Range("F3").Value = 42.36698793397927
$ws.Range("G3").Value = 37.43477470552878
$ws.Range("H3").Value = 16.51701498894585
$ws.Range("I3").Value = 24.21251245331357
$ws.Range("J3").Value = 8.710127221970373
$ws.Range("K3").Value = 8.89696463106198
$ws.Range("M3").Value = 17.35787824235884

$ws.Range("B4").Value = 8.308180316357701
$ws.Range("C4").Value = 4.040950026356565
$ws.Range("E4").Value = 19.01174898619653
$ws.Range("F4").Value = 42.3353150814053
$ws.Range("G4").Value = 37.48559675478093
$ws.Range("H4").Value = 16.55186636563571
$ws.Range("I4").Value = 24.27282016539318
$ws.Range("J4").Value = 8.722618420645389
$ws.Range("K4").Value = 8.800100149713201
$ws.Range("M4").Value = 17.29796885382675

$ws.Range("B5").Value = 8.243474948790711
$ws.Range("C5").Value = 3.99268744588065
$ws.Range("E5").Value = 18.9911465735757
$ws.Range("F5").Value = 42.32507946823695
$ws.Range("G5").Value = 37.50916320492284
$ws.Range("H5").Value = 16.56677052287765
$ws.Range("I5").Value = 24.29857488001153
$ws.Range("J5").Value = 8.727856238970551
$ws.Range("K5").Value = 8.760774960481234
$ws.Range("M5").Value = 17.27431877109274

$ws.Range("B6").Value = 8.232706948305482
$ws.Range("C6").Value = 3.984629349847567
$ws.Range("E6").Value = 18.98778066359114
$ws.Range("F6").Value = 42.32354136551864
$ws.Range("G6").Value = 37.5132484996772
$ws.Range("H6").Value = 16.56928771852964
$ws.Range("I6").Value = 24.30292256400507
$ws.Range("J6").Value = 8.72873489997364
$ws.Range("K6").Value = 8.754255531732667
$ws.Range("M6").Value = 17.27043836135772

$ws.Range("B7").Value = 8.307309337685846
$ws.Range("C7").Value = 4.040302151979929
$ws.Range("E7").Value = 19.01146745020625
$ws.Range("F7").Value = 42.33516621620628
$ws.Range("G7").Value = 37.48590303386487
$ws.Range("H7").Value = 16.55206452728516
$ws.Range("I7").Value = 24.27316273242586
$ws.Range("J7").Value = 8.722688461655283
$ws.Range("K7").Value = 8.799569127365038
$ws.Range("M7").Value = 17.2976467839446

$ws.Range("B8").Value = 8.632323160447889
$ws.Range("C8").Value = 4.279011740352241
$ws.Range("E8").Value = 19.123894976867
$ws.Range("F8").Value = 42.40989669381896
$ws.Range("G8").Value = 37.39034139547813
$ws.Range("H8").Value = 16.48222695747703
$ws.Range("I8").Value = 24.15219428501384
$ws.Range("J8").Value = 8.697301030149855
$ws.Range("K8").Value = 9.000684073334151
$ws.Range("M8").Value = 17.42436150646015

$ws.Range("B9").Value = 9.24809832051926
$ws.Range("C9").Value = 4.717620999484469
$ws.Range("E9").Value = 19.37166413922158
$ws.Range("F9").Value = 42.64006370414654
$ws.Range("G9").Value = 37.29396672966528
$ws.Range("H9").Value = 16.36748586445146
$ws.Range("I9").Value = 23.95234264343255
$ws.Range("J9").Value = 8.652136885422026
$ws.Range("K9").Value = 9.395772389881431
$ws.Range("M9").Value = 17.69539665677314

$ws.Range("B10").Value = 9.68218253400053
$ws.Range("C10").Value = 5.018656595517494
$ws.Range("E10").Value = 19.56925981111982
$ws.Range("F10").Value = 42.85941038002436
$ws.Range("G10").Value = 37.27924585169797
$ws.Range("H10").Value = 16.29677401366723
$ws.Range("I10").Value = 23.82842746332562
$ws.Range("J10").Value = 8.621741849314722
$ws.Range("K10").Value = 9.683484854411361
$ws.Range("M10").Value = 17.9069772169128

$ws.Range("B11").Value = 9.874731432564438
$ws.Range("C11").Value = 5.15046549504719
$ws.Range("E11").Value = 19.66227863253294
$ws.Range("F11").Value = 42.96989952945817
$ws.Range("G11").Value = 37.28482268318681
$ws.Range("H11").Value = 16.26756625947107
$ws.Range("I11").Value = 23.77706037288453
$ws.Range("J11").Value = 8.608513194070628
$ws.Range("K11").Value = 9.813216238427932
$ws.Range("M11").Value = 18.00564143448335

$ws.Range("B12").Value = 9.946868201672016
$ws.Range("C12").Value = 5.19960380516007
$ws.Range("E12").Value = 19.69793023082002
$ws.Range("F12").Value = 43.01325717711268
$ws.Range("G12").Value = 37.28870442889979
$ws.Range("H12").Value = 16.25693247314224
$ws.Range("I12").Value = 23.75833082291482
$ws.Range("J12").Value = 8.603589391433669
$ws.Range("K12").Value = 9.862129334923285
$ws.Range("M12").Value = 18.04332462308174

$ws.Range("B13").Value = 9.931367871656963
$ws.Range("C13").Value = 5.189055978104482
$ws.Range("E13").Value = 19.69023338568822
$ws.Range("F13").Value = 43.0038522001993
$ws.Range("G13").Value = 37.28778964646497
$ws.Range("H13").Value = 16.25920366238594
$ws.Range("I13").Value = 23.76233241037789
$ws.Range("J13").Value = 8.604646019885427
$ws.Range("K13").Value = 9.851605252031469
$ws.Range("M13").Value = 18.03519501002225

$ws.Range("B14").Value = 9.880682123442451
$ws.Range("C14").Value = 5.154523853044624
$ws.Range("E14").Value = 19.66520327935179
$ws.Range("F14").Value = 42.97343630374407
$ws.Range("G14").Value = 37.28510654235924
$ws.Range("H14").Value = 16.26668285699876
$ws.Range("I14").Value = 23.77550499514998
$ws.Range("J14").Value = 8.608106396658043
$ws.Range("K14").Value = 9.817244890648553
$ws.Range("M14").Value = 18.00873536723812

$ws.Range("B15").Value = 9.849532447529329
$ws.Range("C15").Value = 5.133270026700997
$ws.Range("E15").Value = 19.64992658235481
$ws.Range("F15").Value = 42.95500266279725
$ws.Range("G15").Value = 37.28369367869897
$ws.Range("H15").Value = 16.27131965916518
$ws.Range("I15").Value = 23.78366769592521
$ws.Range("J15").Value = 8.610237112892017
$ws.Range("K15").Value = 9.79616900832616
$ws.Range("M15").Value = 17.99256911721268

$ws.Range("B16").Value = 9.66949389315473
$ws.Range("C16").Value = 5.009936010629533
$ws.Range("E16").Value = 19.56324188470501
$ws.Range("F16").Value = 42.85240329397905
$ws.Range("G16").Value = 37.27912883387955
$ws.Range("H16").Value = 16.29874244559878
$ws.Range("I16").Value = 23.83188531788416
$ws.Range("J16").Value = 8.622618374514111
$ws.Range("K16").Value = 9.674979303434593
$ws.Range("M16").Value = 17.900575657621

$ws.Range("B17").Value = 9.557734309641043
$ws.Range("C17").Value = 4.932931850138154
$ws.Range("E17").Value = 19.51084962766754
$ws.Range("F17").Value = 42.79218946049714
$ws.Range("G17").Value = 37.27947608839505
$ws.Range("H17").Value = 16.31632412809496
$ws.Range("I17").Value = 23.86274845806272
$ws.Range("J17").Value = 8.630366792048113
$ws.Range("K17").Value = 9.600303857184635
$ws.Range("M17").Value = 17.84474082450539

$ws.Range("B18").Value = 9.492993966828827
$ws.Range("C18").Value = 4.888159694532227
$ws.Range("E18").Value = 19.48101116710201
$ws.Range("F18").Value = 42.75856533323378
$ws.Range("G18").Value = 37.28083086223437
$ws.Range("H18").Value = 16.3267151428358
$ws.Range("I18").Value = 23.88097082276368
$ws.Range("J18").Value = 8.634879807009542
$ws.Range("K18").Value = 9.557246955270754
$ws.Range("M18").Value = 17.81285503070604

$ws.Range("B19").Value = 9.470997383804461
$ws.Range("C19").Value = 4.872919140016042
$ws.Range("E19").Value = 19.47095993242954
$ws.Range("F19").Value = 42.7473547543721
$ws.Range("G19").Value = 37.28148776387325
$ws.Range("H19").Value = 16.33028116332475
$ws.Range("I19").Value = 23.88722136385795
$ws.Range("J19").Value = 8.636417523886832
$ws.Range("K19").Value = 9.542651996047347
$ws.Range("M19").Value = 17.80209915381503

$ws.Range("B20").Value = 9.569679361987676
$ws.Range("C20").Value = 4.941179180572294
$ws.Range("E20").Value = 19.5163963830574
$ws.Range("F20").Value = 42.79849502053525
$ws.Range("G20").Value = 37.2793195381042
$ws.Range("H20").Value = 16.31442369741056
$ws.Range("I20").Value = 23.85941428856006
$ws.Range("J20").Value = 8.629536132852076
$ws.Range("K20").Value = 9.608264465541074
$ws.Range("M20").Value = 17.850661032356

$ws.Range("B21").Value = 9.89559137677576
$ws.Range("C21").Value = 5.164688069432474
$ws.Range("E21").Value = 19.67254381539417
$ws.Range("F21").Value = 42.98232918292694
$ws.Range("G21").Value = 37.28584656859066
$ws.Range("H21").Value = 16.26447445176589
$ws.Range("I21").Value = 23.77161626999475
$ws.Range("J21").Value = 8.607087680223765
$ws.Range("K21").Value = 9.827343522682879
$ws.Range("M21").Value = 18.01649869919333

$ws.Range("B22").Value = 10.1040352179596
$ws.Range("C22").Value = 5.319703753105265
$ws.Range("E22").Value = 19.77707521802747
$ws.Range("F22").Value = 43.11131059586933
$ws.Range("G22").Value = 37.30042990295473
$ws.Range("H22").Value = 16.23431648017445
$ws.Range("I22").Value = 23.71844510434357
$ws.Range("J22").Value = 8.592915080766407
$ws.Range("K22").Value = 9.969260155373908
$ws.Range("M22").Value = 18.12674291437034

$ws.Range("B23").Value = 9.993223298609335
$ws.Range("C23").Value = 5.231113607555484
$ws.Range("E23").Value = 19.72106567436295
$ws.Range("F23").Value = 43.04167019445078
$ws.Range("G23").Value = 37.2917012416087
$ws.Range("H23").Value = 16.25018449459417
$ws.Range("I23").Value = 23.74643746206043
$ws.Range("J23").Value = 8.600433764493069
$ws.Range("K23").Value = 9.893647475037929
$ws.Range("M23").Value = 18.06774208442597

$ws.Range("B24").Value = 9.564280519271291
$ws.Range("C24").Value = 4.93745212151828
$ws.Range("E24").Value = 19.513887813104
$ws.Range("F24").Value = 42.79564118011304
$ws.Range("G24").Value = 37.27938671662683
$ws.Range("H24").Value = 16.31528200040582
$ws.Range("I24").Value = 23.86092017550553
$ws.Range("J24").Value = 8.629911492090327
$ws.Range("K24").Value = 9.604665860145365
$ws.Range("M24").Value = 17.84798383657331

$ws.Range("B25").Value = 9.084387040249494
$ws.Range("C25").Value = 4.602486642316681
$ws.Range("E25").Value = 19.30182109577991
$ws.Range("F25").Value = 42.5689113008327
$ws.Range("G25").Value = 37.31022647398371
$ws.Range("H25").Value = 16.3961435730848
$ws.Range("I25").Value = 24.00239336721393
$ws.Range("J25").Value = 8.663863408786437
$ws.Range("K25").Value = 9.289129553364972
$ws.Range("M25").Value = 17.6197907153993
